$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was "M", now "B")
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9513888888888888
$ws.Range("C2").Value = 0.958041958041958
$ws.Range("D2").Value = 0.9547038327526133
$ws.Range("E2").Value = 143

# Row 3 (was "B", now "M")
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9285714285714286
$ws.Range("C3").Value = 0.9176470588235294
$ws.Range("D3").Value = 0.9230769230769231
$ws.Range("E3").Value = 85

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9429824561403509
$ws.Range("C4").Value = 0.9429824561403509
$ws.Range("D4").Value = 0.9429824561403509
$ws.Range("E4").Value = 0.9429824561403509

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9399801587301587
$ws.Range("C5").Value = 0.9378445084327437
$ws.Range("D5").Value = 0.9388903779147681
$ws.Range("E5").Value = 228

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9428823795600111
$ws.Range("C6").Value = 0.9429824561403509
$ws.Range("D6").Value = 0.9429130988822902
$ws.Range("E6").Value = 228
